$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 75

# Columns A and D hold text that Excel would otherwise auto-convert
# (a date-looking string and a leading-zero numeric string), so force
# the cell to Text format before assigning the literal string.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-02-18"

$ws.Cells.Item($row, 2).Value = "09:01:16"
$ws.Cells.Item($row, 3).Value = "Tuesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "07"

$ws.Cells.Item($row, 5).Value = 128890
$ws.Cells.Item($row, 6).Value = 140381
$ws.Cells.Item($row, 7).Value = 170910
$ws.Cells.Item($row, 8).Value = 159190
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 145264
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192328
$ws.Cells.Item($row, 14).Value = 115121
$ws.Cells.Item($row, 15).Value = 45405
$ws.Cells.Item($row, 16).Value = 28882
$ws.Cells.Item($row, 17).Value = 66358
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 45823
$ws.Cells.Item($row, 20).Value = -1
